# "feat: cart work mostly"
#
# - Constants sheet: add a new "DISCOUNT CODE" / "welcome15" row (row 8)
# - Address sheet: mark the Massachusetts full_state value as "[D]" (discounted?)
# - Update the remembered selections / active sheet on several worksheets
#
# NOTE: the order of the value writes below matters - it controls the order
# in which brand-new strings are appended to the shared-strings table so the
# resulting <sst> matches what Excel produced (Massachusetts[D], then
# DISCOUNT CODE, then welcome15).

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("Constants")
$wsUser      = $wb.Worksheets.Item("User")
$wsAddress   = $wb.Worksheets.Item("Address")
$wsBvt       = $wb.Worksheets.Item("BVT")

# --- data edits -------------------------------------------------------

# Address!I3 ("full_state" for the shipping row) - "Massachusetts" -> "Massachusetts[D]"
$wsAddress.Range("I3").Value = "Massachusetts[D]"

# Constants: new row 8 - a discount code constant
$wsConstants.Range("A8").Value = "DISCOUNT CODE"
$wsConstants.Range("B8").Value = "welcome15"

# --- selection / active-sheet updates ---------------------------------

$wsUser.Range("A2:G15").Select()

$wsAddress.Range("A2:XFD2").Select()

$wsBvt.Range("S2").Select()

# Constants becomes the active sheet/tab, cursor parked below the new row
$wsConstants.Range("A9").Select()
